# Update the "as_of_utc" timestamp column (AA) on the data sheets.
# The timestamp moves from 2025-11-18 03:06:38 to 2025-11-18 07:06:23
# for rows 2-26 on both the "Главные" and "Линейные" worksheets.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-18 07:06:23"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
